# UserAccount and User API
#
# Adds two new fields to the "User" table documented on Sheet1:
#   - sex       boolean
#   - birthday  datetime
# These are inserted right after "avatar" (row 13) and before "status"
# (previously row 14), so the existing status / mod_user_id / mod_time rows
# move down from rows 14-16 to rows 16-18. Only columns A:C (the User table)
# are affected - columns E:O (other tables documented side-by-side on the
# same sheet) are untouched, so this is done with direct cell writes instead
# of a real row insert.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 17 and 18 previously had no cells in columns A:C - give them the same
# style ("s=2", same as the rest of this table) used throughout A9:C16 by
# copying the format from a row that keeps that style untouched (A9:C9).
$ws.Range("A9:C9").Copy()
$ws.Range("A17:C17").PasteSpecial(-4122)
$ws.Range("A18:C18").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Write the new/shifted values top-to-bottom so the new shared-string
# entries ("sex", "boolean", "birthday") get appended in that order.
$ws.Cells.Item(14, 1).Value = "sex"
$ws.Cells.Item(14, 2).Value = "boolean"
$ws.Cells.Item(14, 3).Value = ""

$ws.Cells.Item(15, 1).Value = "birthday"
$ws.Cells.Item(15, 2).Value = "datetime"
$ws.Cells.Item(15, 3).Value = ""

$ws.Cells.Item(16, 1).Value = "status"
$ws.Cells.Item(16, 2).Value = "int"
$ws.Cells.Item(16, 3).Value = "default: 1"

$ws.Cells.Item(17, 1).Value = "mod_user_id"
$ws.Cells.Item(17, 2).Value = "int"
$ws.Cells.Item(17, 3).Value = ""

$ws.Cells.Item(18, 1).Value = "mod_time"
$ws.Cells.Item(18, 2).Value = "datetime"
$ws.Cells.Item(18, 3).Value = "default: now"

# Match the author's final selection/view state.
$ws.Range("G16").Select()
